$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.62966433333333
$ws.Range("H2").Value = 64.888993
$ws.Range("I2").Value = 0.004276908378962984
$ws.Range("J2").Value = 0.004276908378962984
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.6186053333333333
$ws.Range("N2").Value = 1.855816
$ws.Range("O2").Value = 0.0556943868446899
$ws.Range("P2").Value = 0.0556943868446899
$ws.Range("Q2").Value = 13.38022571480978
$ws.Range("R2").Value = 120.422031433288
$ws.Range("S2").Value = 0.00023819978975726
$ws.Range("T2").Value = 0.00023819978975726

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.62966433333333
$ws.Range("H3").Value = 64.888993
$ws.Range("I3").Value = 0.004276908378962984
$ws.Range("J3").Value = 0.004276908378962984
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.338622
$ws.Range("N3").Value = 7.015866
$ws.Range("O3").Value = 0.2105512373287584
$ws.Range("P3").Value = 0.2105512373287584
$ws.Range("Q3").Value = 50.58360886254867
$ws.Range("R3").Value = 455.252479762938
$ws.Range("S3").Value = 0.0009005083511323906
$ws.Range("T3").Value = 0.0009005083511323906

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.62966433333333
$ws.Range("H4").Value = 64.888993
$ws.Range("I4").Value = 0.004276908378962984
$ws.Range("J4").Value = 0.004276908378962984
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.165314
$ws.Range("N4").Value = 0.495942
$ws.Range("O4").Value = 0.01488357983794147
$ws.Range("P4").Value = 0.01488357983794148
$ws.Range("Q4").Value = 3.575686329600666
$ws.Range("R4").Value = 32.181176966406
$ws.Range("S4").Value = 0.00006365570731785641
$ws.Range("T4").Value = 0.00006365570731785643

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 21.62966433333333
$ws.Range("H5").Value = 64.888993
$ws.Range("I5").Value = 0.004276908378962984
$ws.Range("J5").Value = 0.004276908378962984
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.984598333333333
$ws.Range("N5").Value = 23.953795
$ws.Range("O5").Value = 0.7188707959886103
$ws.Range("P5").Value = 0.7188707959886103
$ws.Range("Q5").Value = 172.7041817864928
$ws.Range("R5").Value = 1554.337636078435
$ws.Range("S5").Value = 0.003074544530755477
$ws.Range("T5").Value = 0.003074544530755477

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4967.017741
$ws.Range("H6").Value = 14901.053223
$ws.Range("I6").Value = 0.9821456064948035
$ws.Range("J6").Value = 0.9821456064948036
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6186053333333333
$ws.Range("N6").Value = 1.855816
$ws.Range("O6").Value = 0.0556943868446899
$ws.Range("P6").Value = 0.0556943868446899
$ws.Range("Q6").Value = 3072.623665343885
$ws.Range("R6").Value = 27653.61298809497
$ws.Range("S6").Value = 0.05469999734593417
$ws.Range("T6").Value = 0.05469999734593418

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 4967.017741
$ws.Range("H7").Value = 14901.053223
$ws.Range("I7").Value = 0.9821456064948035
$ws.Range("J7").Value = 0.9821456064948036
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 2.338622
$ws.Range("N7").Value = 7.015866
$ws.Range("O7").Value = 0.2105512373287584
$ws.Range("P7").Value = 0.2105512373287584
$ws.Range("Q7").Value = 11615.9769634929
$ws.Range("R7").Value = 104543.7926714361
$ws.Range("S7").Value = 0.2067919726844847
$ws.Range("T7").Value = 0.2067919726844848

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4967.017741
$ws.Range("H8").Value = 14901.053223
$ws.Range("I8").Value = 0.9821456064948035
$ws.Range("J8").Value = 0.9821456064948036
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.165314
$ws.Range("N8").Value = 0.495942
$ws.Range("O8").Value = 0.01488357983794147
$ws.Range("P8").Value = 0.01488357983794148
$ws.Range("Q8").Value = 821.1175708356739
$ws.Range("R8").Value = 7390.058137521065
$ws.Range("S8").Value = 0.01461784254674886
$ws.Range("T8").Value = 0.01461784254674886

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4967.017741
$ws.Range("H9").Value = 14901.053223
$ws.Range("I9").Value = 0.9821456064948035
$ws.Range("J9").Value = 0.9821456064948036
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.984598333333333
$ws.Range("N9").Value = 23.953795
$ws.Range("O9").Value = 0.7188707959886103
$ws.Range("P9").Value = 0.7188707959886103
$ws.Range("Q9").Value = 39659.6415764257
$ws.Range("R9").Value = 356936.7741878313
$ws.Range("S9").Value = 0.7060357939176358
$ws.Range("T9").Value = 0.7060357939176359

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.249417666666667
$ws.Range("H10").Value = 6.748253
$ws.Range("I10").Value = 0.0004447851394313067
$ws.Range("J10").Value = 0.0004447851394313068
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.6186053333333333
$ws.Range("N10").Value = 1.855816
$ws.Range("O10").Value = 0.0556943868446899
$ws.Range("P10").Value = 0.0556943868446899
$ws.Range("Q10").Value = 1.391501765494222
$ws.Range("R10").Value = 12.523515889448
$ws.Range("S10").Value = 0.00002477203561825653
$ws.Range("T10").Value = 0.00002477203561825654

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.249417666666667
$ws.Range("H11").Value = 6.748253
$ws.Range("I11").Value = 0.0004447851394313067
$ws.Range("J11").Value = 0.0004447851394313068
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 2.338622
$ws.Range("N11").Value = 7.015866
$ws.Range("O11").Value = 0.2105512373287584
$ws.Range("P11").Value = 0.2105512373287584
$ws.Range("Q11").Value = 5.260537642455333
$ws.Range("R11").Value = 47.344838782098
$ws.Range("S11").Value = 0.00009365006145270596
$ws.Range("T11").Value = 0.00009365006145270597

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.249417666666667
$ws.Range("H12").Value = 6.748253
$ws.Range("I12").Value = 0.0004447851394313067
$ws.Range("J12").Value = 0.0004447851394313068
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.165314
$ws.Range("N12").Value = 0.495942
$ws.Range("O12").Value = 0.01488357983794147
$ws.Range("P12").Value = 0.01488357983794148
$ws.Range("Q12").Value = 0.3718602321473333
$ws.Range("R12").Value = 3.346742089326
$ws.Range("S12").Value = 0.000006619995133455784
$ws.Range("T12").Value = 0.000006619995133455786

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.249417666666667
$ws.Range("H13").Value = 6.748253
$ws.Range("I13").Value = 0.0004447851394313067
$ws.Range("J13").Value = 0.0004447851394313068
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.984598333333333
$ws.Range("N13").Value = 23.953795
$ws.Range("O13").Value = 0.7188707959886103
$ws.Range("P13").Value = 0.7188707959886103
$ws.Range("Q13").Value = 17.96069655223722
$ws.Range("R13").Value = 161.646268970135
$ws.Range("S13").Value = 0.0003197430472268885
$ws.Range("T13").Value = 0.0003197430472268885

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 66.41617433333333
$ws.Range("H14").Value = 199.248523
$ws.Range("I14").Value = 0.01313269998680205
$ws.Range("J14").Value = 0.01313269998680205
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.6186053333333333
$ws.Range("N14").Value = 1.855816
$ws.Range("O14").Value = 0.0556943868446899
$ws.Range("P14").Value = 0.0556943868446899
$ws.Range("Q14").Value = 41.08539966219644
$ws.Range("R14").Value = 369.7685969597679
$ws.Range("S14").Value = 0.0007314176733802075
$ws.Range("T14").Value = 0.0007314176733802075

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 66.41617433333333
$ws.Range("H15").Value = 199.248523
$ws.Range("I15").Value = 0.01313269998680205
$ws.Range("J15").Value = 0.01313269998680205
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 2.338622
$ws.Range("N15").Value = 7.015866
$ws.Range("O15").Value = 0.2105512373287584
$ws.Range("P15").Value = 0.2105512373287584
$ws.Range("Q15").Value = 155.3223264517686
$ws.Range("R15").Value = 1397.900938065918
$ws.Range("S15").Value = 0.002765106231688542
$ws.Range("T15").Value = 0.002765106231688542

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 66.41617433333333
$ws.Range("H16").Value = 199.248523
$ws.Range("I16").Value = 0.01313269998680205
$ws.Range("J16").Value = 0.01313269998680205
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.165314
$ws.Range("N16").Value = 0.495942
$ws.Range("O16").Value = 0.01488357983794147
$ws.Range("P16").Value = 0.01488357983794148
$ws.Range("Q16").Value = 10.97952344374067
$ws.Range("R16").Value = 98.81571099366599
$ws.Range("S16").Value = 0.0001954615887413013
$ws.Range("T16").Value = 0.0001954615887413013

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 66.41617433333333
$ws.Range("H17").Value = 199.248523
$ws.Range("I17").Value = 0.01313269998680205
$ws.Range("J17").Value = 0.01313269998680205
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.984598333333333
$ws.Range("N17").Value = 23.953795
$ws.Range("O17").Value = 0.7188707959886103
$ws.Range("P17").Value = 0.7188707959886103
$ws.Range("Q17").Value = 530.3064748883095
$ws.Range("R17").Value = 4772.758273994785
$ws.Range("S17").Value = 0.009440714492992004
$ws.Range("T17").Value = 0.009440714492992004
